$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.851.67"
$ws.Range("E2").Value = "  +0.40%  "
$ws.Range("D3").Value = "2.810.12"
$ws.Range("E3").Value = "  +1.36%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'356.74"
$ws.Range("D6").Value = "'111.96"
$ws.Range("E6").Value = "  +2.65%  "
$ws.Range("E7").Value = "  +0.86%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "'0.631"
$ws.Range("E9").Value = "  +8.09%  "
$ws.Range("D10").Value = "'40.35"
$ws.Range("E10").Value = "  +2.24%  "
$ws.Range("E11").Value = "  -0.59%  "
$ws.Range("D12").Value = "'0.0840"
$ws.Range("E12").Value = "  -0.40%  "
$ws.Range("D13").Value = "'19.96"
$ws.Range("E13").Value = "  +2.73%  "
$ws.Range("E14").Value = "  +3.15%  "
$ws.Range("D15").Value = "3.253.75"
$ws.Range("E15").Value = "  +1.53%  "
$ws.Range("D16").Value = "2.805.20"
$ws.Range("E16").Value = "  +1.41%  "
$ws.Range("D17").Value = "'0.943"
$ws.Range("E17").Value = "  +1.55%  "
$ws.Range("D18").Value = "51.823.00"
$ws.Range("E18").Value = "  +0.54%  "
$ws.Range("D19").Value = "'7.66"
$ws.Range("E19").Value = "  +2.88%  "
$ws.Range("E20").Value = "  +2.56%  "
$ws.Range("D21").Value = "'13.61"
$ws.Range("E21").Value = "  +4.02%  "
$ws.Range("D22").Value = "0.0₃0978"
$ws.Range("E22").Value = "  +1.17%  "
$ws.Range("D23").Value = "'70.43"
$ws.Range("E23").Value = "  +0.73%  "
$ws.Range("D24").Value = "'268.92"
$ws.Range("E24").Value = "  +0.43%  "
$ws.Range("E25").Value = "  +1.65%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("E27").Value = "  -0.52%  "
$ws.Range("E28").Value = "  -0.84%  "
$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D29").Value = "'38.38"
$ws.Range("E29").Value = "  +11.52%  "
$ws.Range("B30").Value = "Cosmos"
$ws.Range("C30").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D30").Value = "'10.41"
$ws.Range("E30").Value = "  +2.42%  "
$ws.Range("E31").Value = "  +1.56%  "
$ws.Range("D32").Value = "'6.17"
$ws.Range("E32").Value = "  +0.64%  "
$ws.Range("D33").Value = "'52.43"
$ws.Range("E33").Value = "  +1.83%  "
$ws.Range("E34").Value = "  +9.56%  "
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "'0.0881"
$ws.Range("E35").Value = "  +5.42%  "
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").Value = "'0.0444"
$ws.Range("E36").Value = "  -0.45%  "
$ws.Range("D38").Value = "'18.88"
$ws.Range("E38").Value = "  -0.16%  "
$ws.Range("E39").Value = "  +3.24%  "
$ws.Range("D40").Value = "'3.13"
$ws.Range("E40").Value = "  +0.33%  "
$ws.Range("E41").Value = "  +1.37%  "
$ws.Range("E42").Value = "  -1.10%  "
$ws.Range("D43").Value = "'121.03"
$ws.Range("E43").Value = "  +1.32%  "
$ws.Range("D44").Value = "'22.08"
$ws.Range("E44").Value = "  +1.99%  "
$ws.Range("D45").Value = "'2.20"
$ws.Range("E45").Value = "  -0.71%  "
$ws.Range("E46").Value = "  +4.39%  "
$ws.Range("D47").Value = "2.105.68"
$ws.Range("E47").Value = "  +1.04%  "
$ws.Range("D48").Value = "'2.40"
$ws.Range("E48").Value = "  +5.68%  "
$ws.Range("D49").Value = "'0.943"
$ws.Range("E49").Value = "  +0.67%  "
$ws.Range("D50").Value = "'5.45"
$ws.Range("E50").Value = "  -1.62%  "
$ws.Range("E51").Value = "  +7.74%  "
